$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '24.897.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = "'" + '1.711.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.96%  '

$ws.Range("D4").Value = "'" + '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.48%  '

$ws.Range("D5").Value = "'" + '317.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("D6").Value = "'" + '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").Value = "'" + '0.3962'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.44%  '

$ws.Range("D8").Value = "'" + '0.4109'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.15%  '

$ws.Range("D9").Value = "'" + '1.523'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.46%  '

$ws.Range("D10").Value = "'" + '1.005'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.54%  '

$ws.Range("D11").Value = "'" + '53.54'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.23%  '

$ws.Range("D12").Value = "'" + '0.08932'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.96%  '

$ws.Range("D13").Value = "'" + '7.710'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.75%  '

$ws.Range("D14").Value = "'" + '24.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.01%  '

$ws.Range("D15").Value = "'" + '8.169'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").Value = "'" + '0.00001375'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.51%  '

$ws.Range("D17").Value = "'" + '1.695.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.23%  '

$ws.Range("D18").Value = "'" + '100.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("D19").Value = "'" + '0.07138'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.34%  '

$ws.Range("D20").Value = "'" + '20.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.91%  '

$ws.Range("D21").Value = "'" + '7.486'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.63%  '

$ws.Range("D22").Value = "'" + '1.005'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '

$ws.Range("D23").Value = "'" + '14.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.45%  '

$ws.Range("D24").Value = "'" + '24.906.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.99%  '

$ws.Range("D25").Value = "'" + '3.123'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.12%  '

$ws.Range("D26").Value = "'" + '2.332'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("D27").Value = "'" + '23.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '

$ws.Range("D28").Value = "'" + '9.384'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +25.35%  '

$ws.Range("D29").Value = "'" + '165.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.83%  '

$ws.Range("D30").Value = "'" + '139.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.98%  '

$ws.Range("D31").Value = "'" + '5.236'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.41%  '

$ws.Range("D32").Value = "'" + '7.859'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.17%  '

$ws.Range("D33").Value = "'" + '0.09042'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.78%  '

$ws.Range("D34").Value = "'" + '1.883.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("D35").Value = "'" + '1.088'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.62%  '

$ws.Range("D36").Value = "'" + '0.03019'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.91%  '

$ws.Range("D37").Value = "'" + '0.2807'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.02%  '

$ws.Range("D38").Value = "'" + '11.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.94%  '

$ws.Range("D39").Value = "'" + '1.967'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.13%  '

$ws.Range("D40").Value = "'" + '14.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("D41").Value = "'" + '0.09311'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.84%  '

$ws.Range("D42").Value = "'" + '0.8065'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.25%  '

$ws.Range("D43").Value = "'" + '1.486'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.65%  '

$ws.Range("D44").Value = "'" + '16.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.45%  '

$ws.Range("D45").Value = "'" + '0.7359'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.45%  '

$ws.Range("D46").Value = "'" + '2.639'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.43%  '

$ws.Range("D47").Value = "'" + '4.275'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.20%  '

$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = "'" + '1.002'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '

$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").Value = "'" + '1.346'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.33%  '

$ws.Range("D50").Value = "'" + '140.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").Value = "'" + '93.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.26%  '
